# Updates the lattice-multiplication exercise table in place.
# Each cell keeps the same 5-line layout (problem / top digits / ---- / two digit rows)
# but gets new multiplication problems, per the commit diff.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11   # w:br (line break) character used inside a Word Range.Text assignment

$t.Cell(1,1).Range.Text = "89 x 50" + $vt + "  5    0" + $vt + "  ----" + $vt + "8|    |" + $vt + "9|    |"
$t.Cell(1,2).Range.Text = "37 x 21" + $vt + "  2    1" + $vt + "  ----" + $vt + "3|    |" + $vt + "7|    |"
$t.Cell(1,3).Range.Text = "77 x 19" + $vt + "  1    9" + $vt + "  ----" + $vt + "7|    |" + $vt + "7|    |"
$t.Cell(2,1).Range.Text = "47 x 21" + $vt + "  2    1" + $vt + "  ----" + $vt + "4|    |" + $vt + "7|    |"
$t.Cell(2,2).Range.Text = "41 x 64" + $vt + "  6    4" + $vt + "  ----" + $vt + "4|    |" + $vt + "1|    |"
$t.Cell(2,3).Range.Text = "12 x 97" + $vt + "  9    7" + $vt + "  ----" + $vt + "1|    |" + $vt + "2|    |"
$t.Cell(3,1).Range.Text = "25 x 19" + $vt + "  1    9" + $vt + "  ----" + $vt + "2|    |" + $vt + "5|    |"
$t.Cell(3,2).Range.Text = "66 x 15" + $vt + "  1    5" + $vt + "  ----" + $vt + "6|    |" + $vt + "6|    |"
$t.Cell(3,3).Range.Text = "94 x 26" + $vt + "  2    6" + $vt + "  ----" + $vt + "9|    |" + $vt + "4|    |"
$t.Cell(4,1).Range.Text = "23 x 96" + $vt + "  9    6" + $vt + "  ----" + $vt + "2|    |" + $vt + "3|    |"
$t.Cell(4,2).Range.Text = "25 x 52" + $vt + "  5    2" + $vt + "  ----" + $vt + "2|    |" + $vt + "5|    |"
$t.Cell(4,3).Range.Text = "96 x 58" + $vt + "  5    8" + $vt + "  ----" + $vt + "9|    |" + $vt + "6|    |"
$t.Cell(5,1).Range.Text = "58 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "5|    |" + $vt + "8|    |"
$t.Cell(5,2).Range.Text = "20 x 21" + $vt + "  2    1" + $vt + "  ----" + $vt + "2|    |" + $vt + "0|    |"
$t.Cell(5,3).Range.Text = "44 x 58" + $vt + "  5    8" + $vt + "  ----" + $vt + "4|    |" + $vt + "4|    |"
